# Generate Report for Archive
# - "6ee639aa-...md" and "e687fb83-...md" move from "Ready for handoff" to
#   "In Translation" status.
# - The rows describing "c0f34807-...md" and "e687fb83-...md" swap their
#   display order (e687fb83 now sorts before c0f34807).

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, $addr, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 3 (6ee639aa): status -> In Translation
$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"

# Row 4 becomes e687fb83 (was c0f34807), status -> In Translation
$ws.Range("A4").Value = "e687fb83-a44d-4904-b82d-23191b02eef8.md"
$ws.Range("B4").Value = "e2e\e687fb83-a44d-4904-b82d-23191b02eef8.md"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"
$ws.Range("G4").Value = "2016-08-26 22:39:58"

# Row 5 becomes c0f34807 (was e687fb83), status stays Ready for handoff
$ws.Range("A5").Value = "c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md"
$ws.Range("B5").Value = "e2e\c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md"
$ws.Range("G5").Value = "2016-08-26 22:38:51"

Set-HyperlinkDisplay $ws '$B$4' "e2e\e687fb83-a44d-4904-b82d-23191b02eef8.md"
Set-HyperlinkDisplay $ws '$B$5' "e2e\c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 3 (6ee639aa): status -> In Translation
$ws.Range("C3").Value = "In Translation"

# Row 4 becomes e687fb83 (was c0f34807), status -> In Translation
$ws.Range("A4").Value = "e687fb83-a44d-4904-b82d-23191b02eef8.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "e687fb83-a44d-4904-b82d-23191b02eef8.6375e75f55d9e1b6482c2a035ae309da5a7642d5.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-26 22:39:53"

# Row 5 becomes c0f34807 (was e687fb83), status stays Ready for handoff
$ws.Range("A5").Value = "c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md"
$ws.Range("G5").Value = "c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.6072aa6864c71d0aebdc8c31d1d8eed8e0ed6776.zh-cn.xlf"
$ws.Range("H5").Value = "2016-08-26 22:38:47"

Set-HyperlinkDisplay $ws '$A$4' "e687fb83-a44d-4904-b82d-23191b02eef8.md"
Set-HyperlinkDisplay $ws '$A$5' "c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 3 (6ee639aa): status -> In Translation
$ws.Range("C3").Value = "In Translation"

# Row 4 becomes e687fb83 (was c0f34807), status -> In Translation
$ws.Range("A4").Value = "e687fb83-a44d-4904-b82d-23191b02eef8.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "e687fb83-a44d-4904-b82d-23191b02eef8.6375e75f55d9e1b6482c2a035ae309da5a7642d5.de-de.xlf"
$ws.Range("H4").Value = "2016-08-26 22:39:58"

# Row 5 becomes c0f34807 (was e687fb83), status stays Ready for handoff
$ws.Range("A5").Value = "c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md"
$ws.Range("G5").Value = "c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.6072aa6864c71d0aebdc8c31d1d8eed8e0ed6776.de-de.xlf"
$ws.Range("H5").Value = "2016-08-26 22:38:51"

Set-HyperlinkDisplay $ws '$A$4' "e687fb83-a44d-4904-b82d-23191b02eef8.md"
Set-HyperlinkDisplay $ws '$A$5' "c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md"
